# The "Infos à extraire" sheet has two cells (B4, D4) that contained the
# stray placeholder text "n2" / "o2". Clear them back to an (explicitly
# text-typed) empty value, the same way a user would by typing a single
# apostrophe into the cell - this preserves the cell's "Text" formatting
# (quote-prefix) while leaving it blank, instead of fully clearing the
# cell format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "'"
$ws.Range("D4").Value = "'"
